# Reporting log update
# - Rows where "Did Harvest Occur?" (column B) is currently "No":
#     * B -> "Yes"
#     * F (Species, currently blank) -> "Na"
#     * J (Unknown Sex Count, currently 0) -> 1
# - Rows where "Did Harvest Occur?" (column B) is already "Yes" and the
#   Species column (F) was recorded in all caps:
#     * Normalize the species text casing (SQUIRREL -> Squirrel, WEASEL -> Weasel)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToMarkHarvest = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,75,77,78,79)

foreach ($r in $rowsToMarkHarvest) {
    $ws.Range("B$r").Value = "Yes"
    $ws.Range("F$r").Value = "Na"
    $ws.Range("J$r").Value = 1
}

$rowsToFixSpeciesCasing = @(49,50,71,72,73,74,76)

foreach ($r in $rowsToFixSpeciesCasing) {
    $species = $ws.Range("F$r").Value()
    if ($species -eq "SQUIRREL") {
        $ws.Range("F$r").Value = "Squirrel"
    } elseif ($species -eq "WEASEL") {
        $ws.Range("F$r").Value = "Weasel"
    }
}
